# Apply the edits described by the diff:
# 1. Change the shared string "Penyelia Settlement" -> "Penyelia Teller" (cell I2, ROLE column)
# 2. Update the active sheet view: scroll so column D is the top-left visible column,
#    and move the selection from J2 to M2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the ROLE value in cell I2
$ws.Range("I2").Value = "Penyelia Teller"

# 2. Update the window/view: scroll position and selection
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("M2").Select()
